$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is unambiguous (not parseable as a clean number) -
# plain .Value assignment keeps them as text, matching the source diff.
$ws.Range("D2").Value = '42.689.92'
$ws.Range("E2").Value = '  -5.69%  '
$ws.Range("D3").Value = '2.220.45'
$ws.Range("E3").Value = '  -6.15%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("E5").Value = '  +2.02%  '
$ws.Range("E6").Value = '  -10.31%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  -9.28%  '
$ws.Range("E10").Value = '  -11.58%  '
$ws.Range("E11").Value = '  -3.41%  '
$ws.Range("E12").Value = '  -10.58%  '
$ws.Range("E13").Value = '  -9.95%  '
$ws.Range("E14").Value = '  -4.01%  '
$ws.Range("D15").Value = '2.552.92'
$ws.Range("E15").Value = '  -6.14%  '
$ws.Range("E16").Value = '  -13.49%  '
$ws.Range("E17").Value = '  -9.16%  '
$ws.Range("D18").Value = '2.170.95'
$ws.Range("E18").Value = '  -8.06%  '
$ws.Range("D19").Value = '42.528.70'
$ws.Range("E19").Value = '  -5.77%  '
$ws.Range("E20").Value = '  +4.24%  '
$ws.Range("D21").Value = '0.0₃0960'
$ws.Range("E21").Value = '  -10.04%  '
$ws.Range("E22").Value = '  -11.49%  '
$ws.Range("E23").Value = '  -10.81%  '
$ws.Range("E24").Value = '  -7.87%  '
$ws.Range("E25").Value = '  -9.63%  '
$ws.Range("E26").Value = '  -6.80%  '
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("E28").Value = '  -8.81%  '
$ws.Range("E29").Value = '  -8.07%  '
$ws.Range("E30").Value = '  -13.94%  '
$ws.Range("E31").Value = '  -8.83%  '
$ws.Range("E32").Value = '  -8.73%  '
$ws.Range("E33").Value = '  -6.99%  '
$ws.Range("E34").Value = '  -10.77%  '
$ws.Range("E35").Value = '  -6.98%  '
$ws.Range("E36").Value = '  +7.11%  '
$ws.Range("E37").Value = '  -6.94%  '
$ws.Range("E38").Value = '  +8.13%  '
$ws.Range("E39").Value = '  -7.71%  '
$ws.Range("E40").Value = '  -11.62%  '
$ws.Range("E41").Value = '  -10.85%  '
$ws.Range("E42").Value = '  -10.51%  '
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("D44").Value = '1.781.97'
$ws.Range("E44").Value = '  +9.62%  '
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("E45").Value = '  -8.47%  '
$ws.Range("B46").Value = 'BitcoinSV'
$ws.Range("C46").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("E46").Value = '  -13.23%  '
$ws.Range("E47").Value = '  -11.73%  '
$ws.Range("E48").Value = '  -5.82%  '
$ws.Range("E49").Value = '  -4.44%  '
$ws.Range("E50").Value = '  -14.02%  '
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("E51").Value = '  -10.44%  '

# Cells whose new text WOULD be auto-coerced to a number by Excel (e.g. "2.17")
# even though the workbook stores them as plain text ("t=inlineStr", no style).
# Force text entry via NumberFormat "@", then reset the style back to Normal so
# no stray numFmt/style id is left on the cell (matches the unstyled target cells).
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D13",
    "D16",
    "D17",
    "D20",
    "D22",
    "D25",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D41",
    "D45",
    "D46",
    "D47",
    "D48",
    "D50",
    "D51",
)
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }
$ws.Range("D5").Value = '315.50'
$ws.Range("D6").Value = '98.63'
$ws.Range("D7").Value = '0.573'
$ws.Range("D9").Value = '0.559'
$ws.Range("D10").Value = '36.50'
$ws.Range("D13").Value = '7.66'
$ws.Range("D16").Value = '0.855'
$ws.Range("D17").Value = '14.03'
$ws.Range("D20").Value = '13.67'
$ws.Range("D22").Value = '6.48'
$ws.Range("D25").Value = '235.73'
$ws.Range("D28").Value = '10.11'
$ws.Range("D29").Value = '2.17'
$ws.Range("D30").Value = '6.41'
$ws.Range("D31").Value = '20.48'
$ws.Range("D32").Value = '0.0879'
$ws.Range("D33").Value = '158.01'
$ws.Range("D34").Value = '33.82'
$ws.Range("D35").Value = '2.71'
$ws.Range("D36").Value = '3.21'
$ws.Range("D37").Value = '0.122'
$ws.Range("D38").Value = '1.87'
$ws.Range("D39").Value = '4.43'
$ws.Range("D41").Value = '3.51'
$ws.Range("D45").Value = '11.98'
$ws.Range("D46").Value = '87.93'
$ws.Range("D47").Value = '0.206'
$ws.Range("D48").Value = '77.85'
$ws.Range("D50").Value = '60.18'
$ws.Range("D51").Value = '8.43'
foreach ($c in $textCells) { $ws.Range($c).Style = "Normal" }
